$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in header row data (I1/J1/K1) and row 2..4 data in the exact order
#     needed to reproduce the shared-string table order of the target file. ---

$ws.Range("B2").Value = "baseline- full finetune"
$ws.Range("D2").Value = "resnet50"
$ws.Range("I1").Value = "best validation meanAP"
$ws.Range("J1").Value = "test meanAP"
$ws.Range("E2").Value = "SGD"
$ws.Range("F2").Value = "LinearLR"
$ws.Range("C2").Value = "haidar"
$ws.Range("K1").Value = "command"
$ws.Range("K2").Value = 'python src/train.py -m resnet50 -d "./data" -o "./test_model__unfz_lr005" --lr 0.005 -b 2'
$ws.Range("B3").Value = "baseline- full finetune_increaseLR"
$ws.Range("C3").Value = "babar"
$ws.Range("J2").Value = "running"
$ws.Range("K3").Value = 'python src/train.py -m resnet50 -d "./data" -o "./test_model_unfz_lr01" --lr 0.01 -b 2'
$ws.Range("I2").Value = "0.12 at E5"
$ws.Range("B4").Value = "baseline- full finetune_decreaseLR"

# numeric / reused-string cells, row 2
$ws.Range("A2").Value = 1
$ws.Range("G2").Value = 0.005
$ws.Range("H2").Value = 2

# row 3
$ws.Range("A3").Value = 2
$ws.Range("D3").Value = "resnet50"
$ws.Range("E3").Value = "SGD"
$ws.Range("F3").Value = "LinearLR"
$ws.Range("G3").Value = 0.01
$ws.Range("H3").Value = 2
$ws.Range("I3").Value = "running"
$ws.Range("J3").Value = "running"

# row 4
$ws.Range("A4").Value = 3
$ws.Range("D4").Value = "resnet50"
$ws.Range("E4").Value = "SGD"
$ws.Range("F4").Value = "LinearLR"
$ws.Range("G4").Value = 0.001
$ws.Range("H4").Value = 2

# rows 5..19, column A only (id numbers 4..18)
for ($i = 5; $i -le 19; $i++) {
  $ws.Cells.Item($i, 1).Value = $i - 1
}

# --- Styling ---
$ws.Range("A1:K1").Font.Bold = $true

# --- Hyperlink on I2 ---
$ws.Hyperlinks.Add($ws.Range("I2"), "http://example.com", [System.Reflection.Missing]::Value, "0.12@Epoch5")
# restore text since Hyperlinks.Add may touch it
$ws.Range("I2").Value = "0.12 at E5"

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 9.736979166666666
$ws.Columns.Item(2).ColumnWidth = 30.877604166666668
$ws.Columns.Item(3).ColumnWidth = 7.451822916666667
$ws.Columns.Item(4).ColumnWidth = 11.592447916666666
$ws.Columns.Item(5).ColumnWidth = 9.592447916666666
$ws.Columns.Item(6).ColumnWidth = 11.736979166666666
$ws.Columns.Item(7).ColumnWidth = 14.022135416666666
$ws.Columns.Item(8).ColumnWidth = 10.592447916666666
$ws.Columns.Item(9).ColumnWidth = 21.592447916666668
$ws.Columns.Item(10).ColumnWidth = 12.022135416666666
$ws.Columns.Item(11).ColumnWidth = 76.59244791666667

# --- Selection ---
$ws.Range("I4").Select()

Write-Host "done"
